$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force-text cells: values that look like plain decimals would otherwise be
# auto-converted to numbers by Excel, losing the exact "xx.xx0"-style formatting
# and the original text cell type. Setting NumberFormat to "@" (Text) first keeps
# the literal string intact, matching the source inlineStr cells.

$ws.Range("D2").Value = "38.818.77"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.105.21"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.50"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  +2.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.389"
$ws.Range("E9").Value = "  +2.32%  "
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("E11").Value = "  -0.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.84"
$ws.Range("E12").Value = "  +7.01%  "
$ws.Range("D13").Value = "2.418.26"
$ws.Range("E13").Value = "  +1.16%  "
$ws.Range("E14").Value = "  -1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.807"
$ws.Range("E15").Value = "  +2.88%  "
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "2.115.82"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "38.833.43"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.67"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  +1.63%  "
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.07"
$ws.Range("E22").Value = "  +1.24%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.69"
$ws.Range("E26").Value = "  +2.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "171.99"
$ws.Range("E27").Value = "  +1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.137"
$ws.Range("E28").Value = "  +1.39%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.42"
$ws.Range("E29").Value = "  +3.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.34"
$ws.Range("E30").Value = "  +1.76%  "
$ws.Range("E31").Value = "  +9.98%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.59"
$ws.Range("E33").Value = "  +2.14%  "
$ws.Range("E34").Value = "  +12.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.75"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0618"
$ws.Range("E36").Value = "  +2.23%  "
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E38").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.06"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.06"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0227"
$ws.Range("E42").Value = "  +3.49%  "
$ws.Range("D43").Value = "1.525.53"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +8.05%  "
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0915"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.77"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.17"
$ws.Range("E49").Value = "  +1.62%  "
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "2.304.38"
$ws.Range("E51").Value = "  +1.20%  "
